# PrezConfig.xlsx — add the "Excluido" column to the "Recursos" sheet,
# resize columns, normalize page margins, and make "Usuarios" the
# active/selected sheet (it was "Recursos" before).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Usuarios")
$ws2 = $wb.Worksheets.Item("Recursos")

# --- Recursos: new "Excluido" header in column C, matching header style ---
$ws2.Range("C1").Value = "Excluido"
$ws2.Range("B1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 28.8
$ws2.Range("A:B").ColumnWidth = 47.8

# --- Page margins (inches -> points, 72 pt/in) restored to Excel defaults ---
foreach ($ws in @($ws1, $ws2)) {
    $ws.PageSetup.LeftMargin   = 0.7 * 72
    $ws.PageSetup.RightMargin  = 0.7 * 72
    $ws.PageSetup.TopMargin    = 0.75 * 72
    $ws.PageSetup.BottomMargin = 0.75 * 72
    $ws.PageSetup.HeaderMargin = 0.3 * 72
    $ws.PageSetup.FooterMargin = 0.3 * 72
}

# --- Active sheet moves from "Recursos" to "Usuarios" ---
$ws1.Activate()
